$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column H (year 2020), copying column G's formatting (incl.
#     the thick-bordered header-separator row 2) ---
$ws.Range("G2:G5").Copy() | Out-Null
$ws.Range("H2:H5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H3").Value = 2020
$ws.Range("H4").Value = 158
$ws.Range("H5").Value = 397

# --- Update existing 2019 values ---
$ws.Range("G4").Value = 146
$ws.Range("G5").Value = 127

# --- Remove the "Abducted"/"Lost" rows (old rows 6 & 7) ---
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# --- Row 5 (the new last data row) now needs the thick bottom border used
#     by the other boundary rows (row 2 already uses that exact style) ---
$ws.Range("A2:H2").Copy() | Out-Null
$ws.Range("A5:H5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# restore the values/text that the format-only paste did not touch
$ws.Range("A5").Value = "Алынган ок атуучу куралдар"
$ws.Range("B5").Value = "Изъятые огнестрельные оружия"
$ws.Range("C5").Value = "Seized firearms"
$ws.Range("D5").Value = 217
$ws.Range("E5").Value = 399
$ws.Range("F5").Value = 296
$ws.Range("G5").Value = 127
$ws.Range("H5").Value = 397
